$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (column D) cells that are being updated to Text format so that
# numeric-looking strings (e.g. "300.74") are not auto-converted to numbers,
# matching the original inline-string representation in the workbook.
$dCells = @("D2","D3","D5","D7","D8","D9","D11","D12","D13","D14","D15","D16","D18","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D41","D42","D43","D44","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.763.04'
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("D3").Value = '1.867.91'
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '300.74'
$ws.Range("E5").Value = '  -2.21%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '0.5331'
$ws.Range("E7").Value = '  +1.59%  '
$ws.Range("D8").Value = '0.3742'
$ws.Range("E8").Value = '  -1.87%  '
$ws.Range("D9").Value = '0.07182'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("D11").Value = '0.8886'
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").Value = '0.08145'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").Value = '1.880.13'
$ws.Range("E13").Value = '  +20.49%  '
$ws.Range("D14").Value = '92.89'
$ws.Range("E14").Value = '  -3.51%  '
$ws.Range("D15").Value = '5.292'
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = '0.000008492'
$ws.Range("E18").Value = '  -2.14%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '26.792.94'
$ws.Range("D21").Value = '4.982'
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("D23").Value = '6.388'
$ws.Range("E23").Value = '  -1.85%  '
$ws.Range("D24").Value = '2.298'
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").Value = '146.05'
$ws.Range("E25").Value = '  -2.54%  '
$ws.Range("D26").Value = '1.730'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").Value = '18.02'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").Value = '113.91'
$ws.Range("D29").Value = '4.718'
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("D30").Value = '4.622'
$ws.Range("E30").Value = '  -4.86%  '
$ws.Range("D31").Value = '0.09136'
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("D32").Value = '0.8044'
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("D33").Value = '0.05013'
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("E34").Value = '  -4.95%  '
$ws.Range("D35").Value = '2.934'
$ws.Range("E35").Value = '  -1.71%  '
$ws.Range("D36").Value = '0.6135'
$ws.Range("E36").Value = '  +6.40%  '
$ws.Range("D37").Value = '2.663'
$ws.Range("E37").Value = '  -2.18%  '
$ws.Range("D38").Value = '3.194'
$ws.Range("E38").Value = '  -4.88%  '
$ws.Range("E39").Value = '  -2.72%  '
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("D41").Value = '6.508'
$ws.Range("E41").Value = '  -1.70%  '
$ws.Range("D42").Value = '8.763'
$ws.Range("E42").Value = '  -3.81%  '
$ws.Range("D43").Value = '0.5192'
$ws.Range("E43").Value = '  +5.73%  '
$ws.Range("D44").Value = '114.95'
$ws.Range("E44").Value = '  -1.09%  '
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = '1.646'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").Value = '9.905'
$ws.Range("E48").Value = '  -3.10%  '
$ws.Range("D49").Value = '37.63'
$ws.Range("E49").Value = '  -3.08%  '
$ws.Range("D50").Value = '0.06037'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").Value = '62.13'
$ws.Range("E51").Value = '  -3.58%  '

# Reset the style on the Price cells back to Normal so no stray number format
# is left attached to the cell (keeps styles identical to the original).
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
